# Refresh the crypto price/volume snapshot (and the Cronos -> Decentraland
# row-51 swap) to match the latest scrape.
#
# Cells D2:D50/D51 hold price strings that sometimes look like plain numbers
# (e.g. "0.9992"); Excel's COM layer auto-coerces a bare numeric-looking
# .Value assignment into a real number cell. To keep these as text (as the
# source data always stored them), force NumberFormat="@" before the write,
# then restore the "Normal" style afterwards so no stray formatting lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '30.748.04'; Numeric = $false },
    @{ Cell = 'E2'; Value = '  +0.78%  '; Numeric = $false },
    @{ Cell = 'D3'; Value = '1.892.61'; Numeric = $false },
    @{ Cell = 'E3'; Value = '  +1.14%  '; Numeric = $false },
    @{ Cell = 'E4'; Value = '  -0.03%  '; Numeric = $false },
    @{ Cell = 'D5'; Value = '247.13'; Numeric = $true },
    @{ Cell = 'E5'; Value = '  +0.03%  '; Numeric = $false },
    @{ Cell = 'D6'; Value = '0.9992'; Numeric = $true },
    @{ Cell = 'E6'; Value = '  -0.02%  '; Numeric = $false },
    @{ Cell = 'D7'; Value = '0.4728'; Numeric = $true },
    @{ Cell = 'E7'; Value = '  -0.13%  '; Numeric = $false },
    @{ Cell = 'D8'; Value = '0.2925'; Numeric = $true },
    @{ Cell = 'E8'; Value = '  +0.29%  '; Numeric = $false },
    @{ Cell = 'D9'; Value = '0.06520'; Numeric = $true },
    @{ Cell = 'E9'; Value = '  +0.49%  '; Numeric = $false },
    @{ Cell = 'D10'; Value = '22.60'; Numeric = $true },
    @{ Cell = 'E10'; Value = '  +1.88%  '; Numeric = $false },
    @{ Cell = 'D11'; Value = '0.07781'; Numeric = $true },
    @{ Cell = 'E11'; Value = '  +0.87%  '; Numeric = $false },
    @{ Cell = 'D12'; Value = '0.7415'; Numeric = $true },
    @{ Cell = 'E12'; Value = '  -0.25%  '; Numeric = $false },
    @{ Cell = 'D13'; Value = '1.885.79'; Numeric = $false },
    @{ Cell = 'E13'; Value = '  +0.80%  '; Numeric = $false },
    @{ Cell = 'D14'; Value = '96.70'; Numeric = $true },
    @{ Cell = 'E14'; Value = '  -0.70%  '; Numeric = $false },
    @{ Cell = 'D15'; Value = '5.242'; Numeric = $true },
    @{ Cell = 'E15'; Value = '  +1.72%  '; Numeric = $false },
    @{ Cell = 'D16'; Value = '284.77'; Numeric = $true },
    @{ Cell = 'E16'; Value = '  +4.06%  '; Numeric = $false },
    @{ Cell = 'D17'; Value = '30.732.11'; Numeric = $false },
    @{ Cell = 'E17'; Value = '  +0.77%  '; Numeric = $false },
    @{ Cell = 'E18'; Value = '  -1.10%  '; Numeric = $false },
    @{ Cell = 'D19'; Value = '0.000007516'; Numeric = $true },
    @{ Cell = 'E19'; Value = '  +0.07%  '; Numeric = $false },
    @{ Cell = 'D20'; Value = '0.9993'; Numeric = $true },
    @{ Cell = 'E20'; Value = '  -0.09%  '; Numeric = $false },
    @{ Cell = 'D21'; Value = '2.130.57'; Numeric = $false },
    @{ Cell = 'E21'; Value = '  +0.79%  '; Numeric = $false },
    @{ Cell = 'D22'; Value = '5.322'; Numeric = $true },
    @{ Cell = 'E22'; Value = '  +1.37%  '; Numeric = $false },
    @{ Cell = 'D23'; Value = '0.9987'; Numeric = $true },
    @{ Cell = 'E23'; Value = '  -0.03%  '; Numeric = $false },
    @{ Cell = 'D24'; Value = '6.275'; Numeric = $true },
    @{ Cell = 'E24'; Value = '  +1.78%  '; Numeric = $false },
    @{ Cell = 'D25'; Value = '9.226'; Numeric = $true },
    @{ Cell = 'E25'; Value = '  -0.65%  '; Numeric = $false },
    @{ Cell = 'D26'; Value = '164.34'; Numeric = $true },
    @{ Cell = 'E26'; Value = '  +0.42%  '; Numeric = $false },
    @{ Cell = 'D27'; Value = '18.98'; Numeric = $true },
    @{ Cell = 'E27'; Value = '  +1.23%  '; Numeric = $false },
    @{ Cell = 'D28'; Value = '1.922'; Numeric = $true },
    @{ Cell = 'E28'; Value = '  +0.19%  '; Numeric = $false },
    @{ Cell = 'D29'; Value = '1.343'; Numeric = $true },
    @{ Cell = 'E29'; Value = '  -0.42%  '; Numeric = $false },
    @{ Cell = 'D30'; Value = '0.09785'; Numeric = $true },
    @{ Cell = 'E30'; Value = '  -2.10%  '; Numeric = $false },
    @{ Cell = 'D31'; Value = '1.489'; Numeric = $true },
    @{ Cell = 'E31'; Value = '  -1.31%  '; Numeric = $false },
    @{ Cell = 'D32'; Value = '4.299'; Numeric = $true },
    @{ Cell = 'E32'; Value = '  +0.34%  '; Numeric = $false },
    @{ Cell = 'D33'; Value = '4.187'; Numeric = $true },
    @{ Cell = 'E33'; Value = '  +2.09%  '; Numeric = $false },
    @{ Cell = 'D34'; Value = '0.04904'; Numeric = $true },
    @{ Cell = 'E34'; Value = '  +1.96%  '; Numeric = $false },
    @{ Cell = 'D35'; Value = '1.133'; Numeric = $true },
    @{ Cell = 'E35'; Value = '  +1.36%  '; Numeric = $false },
    @{ Cell = 'E36'; Value = '  +0.66%  '; Numeric = $false },
    @{ Cell = 'E37'; Value = '  -0.17%  '; Numeric = $false },
    @{ Cell = 'E38'; Value = '  +2.51%  '; Numeric = $false },
    @{ Cell = 'D39'; Value = '2.837'; Numeric = $true },
    @{ Cell = 'E39'; Value = '  +3.43%  '; Numeric = $false },
    @{ Cell = 'D40'; Value = '76.06'; Numeric = $true },
    @{ Cell = 'E40'; Value = '  +4.29%  '; Numeric = $false },
    @{ Cell = 'D41'; Value = '6.311'; Numeric = $true },
    @{ Cell = 'E41'; Value = '  +1.55%  '; Numeric = $false },
    @{ Cell = 'D42'; Value = '2.014'; Numeric = $true },
    @{ Cell = 'E42'; Value = '  +2.35%  '; Numeric = $false },
    @{ Cell = 'D43'; Value = '0.4297'; Numeric = $true },
    @{ Cell = 'E43'; Value = '  +2.61%  '; Numeric = $false },
    @{ Cell = 'D44'; Value = '0.9993'; Numeric = $true },
    @{ Cell = 'E44'; Value = '  +0.00%  '; Numeric = $false },
    @{ Cell = 'D45'; Value = '0.8333'; Numeric = $true },
    @{ Cell = 'E45'; Value = '  -0.12%  '; Numeric = $false },
    @{ Cell = 'E46'; Value = '  -0.45%  '; Numeric = $false },
    @{ Cell = 'D47'; Value = '9.581'; Numeric = $true },
    @{ Cell = 'E47'; Value = '  +2.50%  '; Numeric = $false },
    @{ Cell = 'D48'; Value = '7.027'; Numeric = $true },
    @{ Cell = 'E48'; Value = '  +0.83%  '; Numeric = $false },
    @{ Cell = 'D49'; Value = '35.42'; Numeric = $true },
    @{ Cell = 'E49'; Value = '  +0.04%  '; Numeric = $false },
    @{ Cell = 'D50'; Value = '910.87'; Numeric = $true },
    @{ Cell = 'E50'; Value = '  -1.40%  '; Numeric = $false },
    @{ Cell = 'B51'; Value = 'Decentraland'; Numeric = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; Numeric = $false },
    @{ Cell = 'D51'; Value = '0.3978'; Numeric = $true },
    @{ Cell = 'E51'; Value = '  +2.82%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
